# Prototypage du système d'interface des registres lié à l'intéractions avec les régiments
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Numeric values (dates / times) - order doesn't matter for these
$ws.Range("A37").Value = 44691
$ws.Range("B37").Value = 0.33333333333333331
$ws.Range("C37").Value = 0.40347222222222223

$ws.Range("A38").Value = 44691
$ws.Range("B38").Value = 0.41666666666666669
$ws.Range("C38").Value = 0.4375

$ws.Range("A39").Value = 44691
$ws.Range("B39").Value = 0.4375
$ws.Range("C39").Value = 0.51041666666666663

$ws.Range("A40").Value = 44691
$ws.Range("B40").Value = 0.5625
$ws.Range("C40").Value = 0.64444444444444449

$ws.Range("A41").Value = 44691
$ws.Range("B41").Value = 0.64444444444444449
$ws.Range("C41").Value = 0.66319444444444442

# Text values - set in the precise order the original author entered them
# so that the shared-strings table is rebuilt with matching indices.
$ws.Range("E37").Value = "Analyse des besoins pour le système régiment"
$ws.Range("F37").Value = "Définition des besoins en terme de comportements`ncomportement similaires entre les composantes`npotentiel candidat pour les interfaces ou abtracts"
$ws.Range("E39").Value = "Prototypage rudimentaire sur base de l'analyse"
$ws.Range("E38").Value = "Préparation de base de l'environnement Unity"
$ws.Range("F38").Value = "Preparation des prefabs notamment les factory précréer afin de faciliter le prototypage"
$ws.Range("F39").Value = "Gros Soucis de conception, il faut reprendre depuis le début et bien cerner les responsabilité de chaque classe"
$ws.Range("E40").Value = "Suite Prototypage"
$ws.Range("F40").Value = "redéfinition des interface selon les observation précédentes"
$ws.Range("G40").Value = "L'implémentation a pris plus de temps, car il a fallu adapter plus d'élément que prévu`nCependant l'architecture smelbe être solide en l'état"
$ws.Range("E41").Value = "Soucis liées aux interface sur Unity"

# Update view to reflect scroll position / selection
$ws.Application.ActiveWindow.ScrollRow = 36
$ws.Range("F41").Select()
